$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the specific cell values that differ between each pair of rows.
# (Columns whose values are identical between the paired rows are left
# untouched so their original types/formatting are preserved.)
function Swap-Cells($ws, $r1, $r2, $cols) {
    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"
        $v1 = $ws.Range($addr1).Value()
        $v2 = $ws.Range($addr2).Value()
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

Swap-Cells $ws 3 4 @("A","B","E","F","G","H","Q","R","AC")
Swap-Cells $ws 14 15 @("A","Q","R")
Swap-Cells $ws 16 17 @("A","Q","R","AC")
